$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell-by-cell updates per the source diff.
# Column D values are forced to Text (NumberFormat "@") before assignment so
# numeric-looking strings (e.g. "1.001", "306.62") are NOT auto-converted to
# numbers by Excel, then the style is reset to Normal so no stray number-format
# style is left attached to the cell (matching the original unstyled cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.086.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5187"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3753"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07216"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9014"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.949.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.57%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07662"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.239"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008505"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.142.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.062"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.156.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.381"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.305"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.726"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.930"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.796"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09213"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05051"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.246"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7713"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.982"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.287"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.583"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5624"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01989"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.028"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.639"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1511"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4846"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.596"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.49%  "
